$wb = $excel.ActiveWorkbook

# ===== Sheet1: Overview =====
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A4").Value = '93931646-dee5-40ed-b807-0f710ced4785.md'
$ws1.Range("B4").Value = 'e2e\93931646-dee5-40ed-b807-0f710ced4785.md'
$ws1.Range("C4").Value = '.md'
$ws1.Range("D4").Value = ''''
$ws1.Range("E4").Value = 'Ready for handoff'
$ws1.Range("F4").Value = 'Ready for handoff'
$ws1.Range("G4").Value = '2016-08-29 16:44:52'

$ws1.Range("A5").Value = 'ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md'
$ws1.Range("B5").Value = 'e2e\ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md'
$ws1.Range("C5").Value = '.md'
$ws1.Range("D5").Value = ''''
$ws1.Range("E5").Value = 'Ready for handoff'
$ws1.Range("F5").Value = 'Ready for handoff'
$ws1.Range("G5").Value = '2016-08-29 16:44:52'

$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fbab17eca8384c18d95238340e13f0351d924b59/e2e/93931646-dee5-40ed-b807-0f710ced4785.md", "", "", "e2e\93931646-dee5-40ed-b807-0f710ced4785.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fbab17eca8384c18d95238340e13f0351d924b59/e2e/ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md", "", "", "e2e\ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md") | Out-Null

$ws1.ListObjects.Item(1).Resize($ws1.Range("A1:G5"))

# ===== Sheet2: zh-cn =====
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A4").Value = '93931646-dee5-40ed-b807-0f710ced4785.md'
$ws2.Range("B4").Value = '.md'
$ws2.Range("C4").Value = 'Ready for handoff'
$ws2.Range("D4").Value = 'e2e'
$ws2.Range("E4").Value = 'ht'
$ws2.Range("F4").Value = '''False'
$ws2.Range("G4").Value = '93931646-dee5-40ed-b807-0f710ced4785.c5a62ea07dc5f11e8bc0e19c5957b7b3015d5501.zh-cn.xlf'
$ws2.Range("H4").Value = '2016-08-29 16:44:47'
$ws2.Range("I4").Value = ''''
$ws2.Range("J4").Value = ''''
$ws2.Range("K4").Value = '0001-01-01 00:00:00'
$ws2.Range("L4").Value = ''''
$ws2.Range("M4").Value = '''True'
$ws2.Range("N4").Value = ''''
$ws2.Range("O4").Value = '''False'
$ws2.Range("P4").Value = ''''

$ws2.Range("A5").Value = 'ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md'
$ws2.Range("B5").Value = '.md'
$ws2.Range("C5").Value = 'Ready for handoff'
$ws2.Range("D5").Value = 'e2e'
$ws2.Range("E5").Value = 'ht'
$ws2.Range("F5").Value = '''False'
$ws2.Range("G5").Value = 'ed420cf4-c719-4f89-a9be-3e4dff3ae91c.9ad492e7902f6a345ce3a361a40e3b0d4a7b2226.zh-cn.xlf'
$ws2.Range("H5").Value = '2016-08-29 16:44:47'
$ws2.Range("I5").Value = ''''
$ws2.Range("J5").Value = ''''
$ws2.Range("K5").Value = '0001-01-01 00:00:00'
$ws2.Range("L5").Value = ''''
$ws2.Range("M5").Value = '''True'
$ws2.Range("N5").Value = ''''
$ws2.Range("O5").Value = '''False'
$ws2.Range("P5").Value = ''''

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fbab17eca8384c18d95238340e13f0351d924b59/e2e/93931646-dee5-40ed-b807-0f710ced4785.md", "", "", "93931646-dee5-40ed-b807-0f710ced4785.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fbab17eca8384c18d95238340e13f0351d924b59/e2e/ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md", "", "", "ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md") | Out-Null

$ws2.ListObjects.Item(1).Resize($ws2.Range("A1:P5"))

# ===== Sheet3: de-de =====
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A4").Value = '93931646-dee5-40ed-b807-0f710ced4785.md'
$ws3.Range("B4").Value = '.md'
$ws3.Range("C4").Value = 'Ready for handoff'
$ws3.Range("D4").Value = 'e2e'
$ws3.Range("E4").Value = 'ht'
$ws3.Range("F4").Value = '''False'
$ws3.Range("G4").Value = '93931646-dee5-40ed-b807-0f710ced4785.c5a62ea07dc5f11e8bc0e19c5957b7b3015d5501.de-de.xlf'
$ws3.Range("H4").Value = '2016-08-29 16:44:52'
$ws3.Range("I4").Value = ''''
$ws3.Range("J4").Value = ''''
$ws3.Range("K4").Value = '0001-01-01 00:00:00'
$ws3.Range("L4").Value = ''''
$ws3.Range("M4").Value = '''True'
$ws3.Range("N4").Value = ''''
$ws3.Range("O4").Value = '''False'
$ws3.Range("P4").Value = ''''

$ws3.Range("A5").Value = 'ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md'
$ws3.Range("B5").Value = '.md'
$ws3.Range("C5").Value = 'Ready for handoff'
$ws3.Range("D5").Value = 'e2e'
$ws3.Range("E5").Value = 'ht'
$ws3.Range("F5").Value = '''False'
$ws3.Range("G5").Value = 'ed420cf4-c719-4f89-a9be-3e4dff3ae91c.9ad492e7902f6a345ce3a361a40e3b0d4a7b2226.de-de.xlf'
$ws3.Range("H5").Value = '2016-08-29 16:44:52'
$ws3.Range("I5").Value = ''''
$ws3.Range("J5").Value = ''''
$ws3.Range("K5").Value = '0001-01-01 00:00:00'
$ws3.Range("L5").Value = ''''
$ws3.Range("M5").Value = '''True'
$ws3.Range("N5").Value = ''''
$ws3.Range("O5").Value = '''False'
$ws3.Range("P5").Value = ''''

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fbab17eca8384c18d95238340e13f0351d924b59/e2e/93931646-dee5-40ed-b807-0f710ced4785.md", "", "", "93931646-dee5-40ed-b807-0f710ced4785.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fbab17eca8384c18d95238340e13f0351d924b59/e2e/ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md", "", "", "ed420cf4-c719-4f89-a9be-3e4dff3ae91c.md") | Out-Null

$ws3.ListObjects.Item(1).Resize($ws3.Range("A1:P5"))

